# Applies updated cryptocurrency price/volume data to sheet1 of the workbook
# (mirrors the automated "Updated cryptos list" GitHub Actions commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.923.80"
$ws.Range("E2").Value = "'  -1.88%  "
$ws.Range("D3").Value = "'1.831.75"
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "'  +0.39%  "
$ws.Range("D5").Value = "'311.08"
$ws.Range("E5").Value = "'  -1.53%  "
$ws.Range("D6").Value = "'1.007"
$ws.Range("E6").Value = "'  +0.46%  "
$ws.Range("D7").Value = "'0.4600"
$ws.Range("E7").Value = "'  -1.55%  "
$ws.Range("D8").Value = "'0.3651"
$ws.Range("E8").Value = "'  -2.73%  "
$ws.Range("D9").Value = "'0.07198"
$ws.Range("E9").Value = "'  -2.77%  "
$ws.Range("D10").Value = "'0.8786"
$ws.Range("E10").Value = "'  -1.48%  "
$ws.Range("D11").Value = "'0.07867"
$ws.Range("E11").Value = "'  -1.24%  "
$ws.Range("E12").Value = "'  -2.61%  "
$ws.Range("D13").Value = "'1.864.34"
$ws.Range("E13").Value = "'  +0.03%  "
$ws.Range("D14").Value = "'5.320"
$ws.Range("E14").Value = "'  -2.40%  "
$ws.Range("D15").Value = "'6.424"
$ws.Range("E15").Value = "'  -3.23%  "
$ws.Range("D16").Value = "'89.72"
$ws.Range("E16").Value = "'  -3.58%  "
$ws.Range("D17").Value = "'1.008"
$ws.Range("E17").Value = "'  +0.45%  "
$ws.Range("D18").Value = "'0.000008730"
$ws.Range("E18").Value = "'  -2.72%  "
$ws.Range("D19").Value = "'1.007"
$ws.Range("E19").Value = "'  +0.47%  "
$ws.Range("D20").Value = "'26.950.60"
$ws.Range("E20").Value = "'  -1.87%  "
$ws.Range("D21").Value = "'14.50"
$ws.Range("E21").Value = "'  -3.19%  "
$ws.Range("D22").Value = "'5.009"
$ws.Range("E22").Value = "'  -3.52%  "
$ws.Range("E23").Value = "'  -1.70%  "
$ws.Range("D24").Value = "'2.002"
$ws.Range("E24").Value = "'  +6.38%  "
$ws.Range("D25").Value = "'150.44"
$ws.Range("E25").Value = "'  -1.70%  "
$ws.Range("D26").Value = "'18.22"
$ws.Range("E26").Value = "'  -2.13%  "
$ws.Range("D27").Value = "'1.997"
$ws.Range("E27").Value = "'  -5.36%  "
$ws.Range("D28").Value = "'114.32"
$ws.Range("E28").Value = "'  -2.76%  "
$ws.Range("D29").Value = "'4.946"
$ws.Range("E29").Value = "'  -4.56%  "
$ws.Range("D30").Value = "'0.08819"
$ws.Range("E30").Value = "'  -1.27%  "
$ws.Range("D31").Value = "'3.107"
$ws.Range("E31").Value = "'  +4.17%  "
$ws.Range("D32").Value = "'0.7626"
$ws.Range("E32").Value = "'  +0.51%  "
$ws.Range("D33").Value = "'4.460"
$ws.Range("E33").Value = "'  -1.39%  "
$ws.Range("D34").Value = "'1.136"
$ws.Range("E34").Value = "'  -2.34%  "
$ws.Range("D35").Value = "'2.664"
$ws.Range("E35").Value = "'  +0.25%  "
$ws.Range("E36").Value = "'  +0.63%  "
$ws.Range("D37").Value = "'0.01929"
$ws.Range("E37").Value = "'  -1.84%  "
$ws.Range("D38").Value = "'0.05153"
$ws.Range("E38").Value = "'  -2.99%  "
$ws.Range("E39").Value = "'  -2.34%  "
$ws.Range("D40").Value = "'6.944"
$ws.Range("E40").Value = "'  -3.95%  "
$ws.Range("D41").Value = "'0.5003"
$ws.Range("E41").Value = "'  -4.84%  "
$ws.Range("D42").Value = "'0.1599"
$ws.Range("E42").Value = "'  -3.15%  "
$ws.Range("D43").Value = "'8.352"
$ws.Range("E43").Value = "'  +0.01%  "
$ws.Range("B44").Value = "'EnergySwap"
$ws.Range("C44").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'10.25"
$ws.Range("E44").Value = "'  -0.99%  "
$ws.Range("D45").Value = "'0.4670"
$ws.Range("E45").Value = "'  -5.23%  "
$ws.Range("B46").Value = "'PaxDollar"
$ws.Range("C46").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.008"
$ws.Range("E46").Value = "'  +0.55%  "
$ws.Range("D47").Value = "'102.76"
$ws.Range("E47").Value = "'  -1.24%  "
$ws.Range("D48").Value = "'1.605"
$ws.Range("E48").Value = "'  -3.72%  "
$ws.Range("D49").Value = "'0.06119"
$ws.Range("E49").Value = "'  -2.51%  "
$ws.Range("D50").Value = "'64.80"
$ws.Range("E50").Value = "'  -2.07%  "
$ws.Range("D51").Value = "'36.29"
$ws.Range("E51").Value = "'  -2.64%  "
